$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update column N values (rows 4-19), including adding a new value at N19
$ws.Range("N4").Value = 9
$ws.Range("N5").Value = 8
$ws.Range("N6").Value = 2
$ws.Range("N7").Value = 2
$ws.Range("N8").Value = 17
$ws.Range("N9").Value = 16
$ws.Range("N10").Value = 27
$ws.Range("N11").Value = 30
$ws.Range("N12").Value = 40
$ws.Range("N13").Value = 32
$ws.Range("N14").Value = 40
$ws.Range("N15").Value = 45
$ws.Range("N16").Value = 52
$ws.Range("N17").Value = 52
$ws.Range("N18").Value = 52
$ws.Range("N19").Value = 52

# Update the active selection from K20 to N20
$ws.Range("N20").Select()
